# Add extra columns (WIN, TOP4, TOP5, RELEGATION) before ExpPoints, which
# moves from C to G. Also refresh the standings (team order + ExpPoints
# values) for matchday 7. Prep work for a Monte Carlo simulation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# C1 used to hold "ExpPoints"; it now becomes "WIN", four new headers are
# inserted after it, and "ExpPoints" is pushed out to G1. Clone the header
# style (bold / bordered / centered) from C1 onto every new header cell
# before overwriting the text.
$ws.Range("C1").Copy()
$ws.Range("D1:G1").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("G1").Value = $ws.Range("C1").Value2
$ws.Range("C1").Value = "WIN"
$ws.Range("D1").Value = "TOP4"
$ws.Range("E1").Value = "TOP5"
$ws.Range("F1").Value = "RELEGATION"

# --- Updated matchday-7 standings ---------------------------------------
# rank, team, expected points (ExpPoints now lives in column G)
$teams = @(
    @(1,  "Arsenal",                  79.66402555895789),
    @(2,  "Manchester City",          72.33285879966623),
    @(3,  "Liverpool",                72.12521128168515),
    @(4,  "Chelsea",                  61.41638994325847),
    @(5,  "Crystal Palace",           59.34894389596209),
    @(6,  "Aston Villa",              56.33905967742943),
    @(7,  "AFC Bournemouth",          55.77841714752951),
    @(8,  "Newcastle United",         55.59877853746523),
    @(9,  "Brighton & Hove Albion",   54.06223694944708),
    @(10, "Tottenham Hotspur",        53.72158121697367),
    @(11, "Manchester United",        50.59775186775695),
    @(12, "Brentford",                48.96878462918131),
    @(13, "Everton",                  45.88313310458392),
    @(14, "Fulham",                   43.71636376435959),
    @(15, "Sunderland",               42.34470993233224),
    @(16, "Nottingham Forest",        38.29705922271609),
    @(17, "Leeds United",             37.00477809322039),
    @(18, "West Ham United",          35.76429468009974),
    @(19, "Burnley",                  32.68901559402691),
    @(20, "Wolverhampton Wanderers",  29.82627784039368)
)

for ($i = 0; $i -lt $teams.Count; $i++) {
    $row = $i + 2
    $rank = $teams[$i][0]
    $name = $teams[$i][1]
    $points = $teams[$i][2]

    $ws.Cells.Item($row, 1).Value = $rank
    $ws.Cells.Item($row, 2).Value = $name

    # New placeholder columns for the upcoming Monte Carlo simulation
    # (win / top4 / top5 / relegation probabilities) - empty for now.
    $ws.Cells.Item($row, 3).Value = ""
    $ws.Cells.Item($row, 4).Value = ""
    $ws.Cells.Item($row, 5).Value = ""
    $ws.Cells.Item($row, 6).Value = ""

    $ws.Cells.Item($row, 7).Value = $points
}
